# Change the table style of the table on slide 6 from the custom
# "Table_0" style ({6462257C-74E7-45B3-8B3C-67C4FA5EC96B}, defined in
# ppt/tableStyles.xml) to the built-in table style
# {385CB8CB-A4D6-479C-8726-88C38E6464E2}.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

foreach ($shp in $s.Shapes) {
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{385CB8CB-A4D6-479C-8726-88C38E6464E2}")
    }
}
